# Apply cryptos list update (GitHub Actions scrape refresh)
# Uses a leading apostrophe to force text entry (prevents Excel from
# auto-converting numeric-looking strings like "215.12" into numbers),
# then resets the cell Style back to "Normal" so no stray number-format
# / quote-prefix style is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.945.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.07%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.654.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.94%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.10%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'215.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.33%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +2.16%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.07%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +2.46%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +1.40%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'20.22"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +4.82%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0874"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +2.02%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.887.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +2.96%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.665.87"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +3.19%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +1.95%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +2.29%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +2.60%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'26.953.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +2.10%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'236.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.52%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +0.85%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'7.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.53%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +0.00%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +3.38%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'9.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +3.82%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +3.29%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'145.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = "'  +1.70%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +0.75%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.01%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'15.79"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +2.18%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.0496"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.33%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +1.32%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.548.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +4.09%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +2.86%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +4.73%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +8.10%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -0.24%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.584"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +3.93%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.893"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +8.62%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.0170"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +2.92%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'6.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +3.24%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +0.04%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +2.88%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'65.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +7.72%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.795.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +2.80%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +1.80%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.916"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.85%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'90.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.10%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +1.49%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +4.30%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0988"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +2.87%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +0.79%  "
$ws.Range("E51").Style = "Normal"
